$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037091173169014
$ws.Range("D2").Value = 1.040282722872028
$ws.Range("E2").Value = 1.04533018781891
$ws.Range("F2").Value = 1.054812182243928
$ws.Range("I2").Value = 1.039464879571445
$ws.Range("J2").Value = 1.04219593268829
$ws.Range("K2").Value = 1.043065512194089
$ws.Range("L2").Value = 1.048098730908623
$ws.Range("M2").Value = 1.057554362391574
$ws.Range("N2").Value = 1.043675970454252
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038009959035218
$ws.Range("D3").Value = 1.040814915377988
$ws.Range("E3").Value = 1.046155947165207
$ws.Range("F3").Value = 1.055741319233102
$ws.Range("I3").Value = 1.039647222211805
$ws.Range("J3").Value = 1.042759108239194
$ws.Range("K3").Value = 1.043408828738656
$ws.Range("L3").Value = 1.048735871530986
$ws.Range("M3").Value = 1.058296520902494
$ws.Range("N3").Value = 1.044239945779036
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038605063990815
$ws.Range("D4").Value = 1.041159535521071
$ws.Range("E4").Value = 1.046691178926865
$ws.Range("F4").Value = 1.056343594624388
$ws.Range("I4").Value = 1.03976425086455
$ws.Range("J4").Value = 1.043123477507471
$ws.Range("K4").Value = 1.043630521926876
$ws.Range("L4").Value = 1.049148404030571
$ws.Range("M4").Value = 1.058777178172012
$ws.Range("N4").Value = 1.044604832493466
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038855385940638
$ws.Range("D5").Value = 1.04130447331623
$ws.Range("E5").Value = 1.046916406758824
$ws.Range("F5").Value = 1.056597043780585
$ws.Range("I5").Value = 1.039813219393678
$ws.Range("J5").Value = 1.043276646987781
$ws.Range("K5").Value = 1.043723611451536
$ws.Range("L5").Value = 1.049321893712781
$ws.Range("M5").Value = 1.058979348333533
$ws.Range("N5").Value = 1.044758219492007
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038897424262311
$ws.Range("D6").Value = 1.041328812445678
$ws.Range("E6").Value = 1.046954236148581
$ws.Range("F6").Value = 1.056639613774929
$ws.Range("I6").Value = 1.039821427899717
$ws.Range("J6").Value = 1.043302364120199
$ws.Range("K6").Value = 1.043739235089089
$ws.Range("L6").Value = 1.049351026937047
$ws.Range("M6").Value = 1.059013299526088
$ws.Range("N6").Value = 1.044783973145704
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038608408257025
$ws.Range("D7").Value = 1.041161471955514
$ws.Range("E7").Value = 1.046694187583777
$ws.Range("F7").Value = 1.056346980233916
$ws.Range("I7").Value = 1.039764906090261
$ws.Range("J7").Value = 1.043125524211542
$ws.Range("K7").Value = 1.043631766228574
$ws.Range("L7").Value = 1.049150721969939
$ws.Range("M7").Value = 1.058779879181121
$ws.Range("N7").Value = 1.044606882104093
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037401558464186
$ws.Range("D8").Value = 1.040462525732061
$ws.Range("E8").Value = 1.045609067671746
$ws.Range("F8").Value = 1.05512596767479
$ws.Range("I8").Value = 1.039526701328732
$ws.Range("J8").Value = 1.042386268754705
$ws.Range("K8").Value = 1.043181631220193
$ws.Range("L8").Value = 1.048314001023999
$ws.Range("M8").Value = 1.057805088259098
$ws.Range("N8").Value = 1.043866576819712
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035279499406879
$ws.Range("D9").Value = 1.039232933377744
$ws.Range("E9").Value = 1.043703986170031
$ws.Range("F9").Value = 1.052982588009006
$ws.Range("I9").Value = 1.039099637446509
$ws.Range("J9").Value = 1.041083324502816
$ws.Range("K9").Value = 1.042385000853009
$ws.Range("L9").Value = 1.046841635677187
$ws.Range("M9").Value = 1.056090743756493
$ws.Range("N9").Value = 1.042561782237536
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033867932828966
$ws.Range("D10").Value = 1.038414684229449
$ws.Range("E10").Value = 1.042438754081929
$ws.Range("F10").Value = 1.051559275235218
$ws.Range("I10").Value = 1.038810049107511
$ws.Range("J10").Value = 1.040214569169955
$ws.Range("K10").Value = 1.041851677548218
$ws.Range("L10").Value = 1.045861504389048
$ws.Range("M10").Value = 1.054950188164211
$ws.Range("N10").Value = 1.041691793172456
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033257468159535
$ws.Range("D11").Value = 1.038060746854566
$ws.Range("E11").Value = 1.041892057267286
$ws.Range("F11").Value = 1.050944314571989
$ws.Range("I11").Value = 1.038683505350477
$ws.Range("J11").Value = 1.039838371946892
$ws.Range("K11").Value = 1.041620227652564
$ws.Range("L11").Value = 1.045437454546074
$ws.Range("M11").Value = 1.054456887055004
$ws.Range("N11").Value = 1.041315061706185
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033030828667016
$ws.Range("D12").Value = 1.037929336172619
$ws.Range("E12").Value = 1.041689165072195
$ws.Range("F12").Value = 1.050716094160791
$ws.Range("I12").Value = 1.038636329238666
$ws.Range("J12").Value = 1.039698633527344
$ws.Range("K12").Value = 1.041534180363718
$ws.Range("L12").Value = 1.045279997847269
$ws.Range("M12").Value = 1.054273739584705
$ws.Range("N12").Value = 1.041175124842054
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033079438418324
$ws.Range("D13").Value = 1.03795752158214
$ws.Range("E13").Value = 1.041732678177808
$ws.Range("F13").Value = 1.050765038986063
$ws.Range("I13").Value = 1.038646456453206
$ws.Range("J13").Value = 1.039728607967273
$ws.Range("K13").Value = 1.041552641240836
$ws.Range("L13").Value = 1.04531377036766
$ws.Range("M13").Value = 1.05431302140174
$ws.Range("N13").Value = 1.041205141849124
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033238731730759
$ws.Range("D14").Value = 1.038049883224647
$ws.Range("E14").Value = 1.041875282549438
$ws.Range("F14").Value = 1.050925445644055
$ws.Range("I14").Value = 1.03867960926584
$ws.Range("J14").Value = 1.039826821160689
$ws.Range("K14").Value = 1.041613116510087
$ws.Range("L14").Value = 1.045424438011778
$ws.Range("M14").Value = 1.054441746255103
$ws.Range("N14").Value = 1.041303494516541
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033336892773559
$ws.Range("D15").Value = 1.038106797944435
$ws.Range("E15").Value = 1.041963169089334
$ws.Range("F15").Value = 1.05102430447421
$ws.Range("I15").Value = 1.038700013019357
$ws.Range("J15").Value = 1.03988733331424
$ws.Range("K15").Value = 1.041650367220463
$ws.Range("L15").Value = 1.04549263122983
$ws.Range("M15").Value = 1.054521069391424
$ws.Range("N15").Value = 1.04136409260429
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033908463320102
$ws.Range("D16").Value = 1.038438181838373
$ws.Range("E16").Value = 1.042475061055411
$ws.Range("F16").Value = 1.05160011658154
$ws.Range("I16").Value = 1.038818423225254
$ws.Range("J16").Value = 1.040239535799697
$ws.Range("K16").Value = 1.041867027320624
$ws.Range("L16").Value = 1.045889654729222
$ws.Range("M16").Value = 1.054982938994507
$ws.Range("N16").Value = 1.041716795257675
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034267196819336
$ws.Range("D17").Value = 1.038646150537748
$ws.Range("E17").Value = 1.042796468199007
$ws.Range("F17").Value = 1.05196166901034
$ws.Range("I17").Value = 1.038892391373295
$ws.Range("J17").Value = 1.040460458383601
$ws.Range("K17").Value = 1.042002794833926
$ws.Range("L17").Value = 1.04613879255431
$ws.Range("M17").Value = 1.055272810339957
$ws.Range("N17").Value = 1.041938031576988
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034476512574382
$ws.Range("D18").Value = 1.038767490728736
$ws.Range("E18").Value = 1.042984051083064
$ws.Range("F18").Value = 1.05217268591544
$ws.Range("I18").Value = 1.038935424730568
$ws.Range("J18").Value = 1.040589316747336
$ws.Range("K18").Value = 1.042081935724272
$ws.Range("L18").Value = 1.046284144457806
$ws.Range("M18").Value = 1.055441942007582
$ws.Range("N18").Value = 1.042067072934381
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034547896112874
$ws.Range("D19").Value = 1.038808870574786
$ws.Range("E19").Value = 1.043048030871202
$ws.Range("F19").Value = 1.052244659143132
$ws.Range("I19").Value = 1.038950079143609
$ws.Range("J19").Value = 1.04063325374118
$ws.Range("K19").Value = 1.04210891220357
$ws.Range("L19").Value = 1.046333711435803
$ws.Range("M19").Value = 1.055499620758519
$ws.Range("N19").Value = 1.042111072323795
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034228700594938
$ws.Range("D20").Value = 1.038623833772311
$ws.Range("E20").Value = 1.042761972706749
$ws.Range("F20").Value = 1.05192286445132
$ws.Range("I20").Value = 1.038884466768751
$ws.Range("J20").Value = 1.040436755684743
$ws.Range("K20").Value = 1.041988233427759
$ws.Range("L20").Value = 1.04611205889552
$ws.Range("M20").Value = 1.05524170422901
$ws.Range("N20").Value = 1.041914295217579
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03319182065702
$ws.Range("D21").Value = 1.038022683423178
$ws.Range("E21").Value = 1.041833284263818
$ws.Range("F21").Value = 1.050878204254984
$ws.Range("I21").Value = 1.038669851335025
$ws.Range("J21").Value = 1.039797899870478
$ws.Range("K21").Value = 1.041595310152221
$ws.Range("L21").Value = 1.045391847640581
$ws.Range("M21").Value = 1.054403837594843
$ws.Range("N21").Value = 1.04127453215478
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032540554766938
$ws.Range("D22").Value = 1.037645049451531
$ws.Range("E22").Value = 1.041250396213603
$ws.Range("F22").Value = 1.050222562283568
$ws.Range("I22").Value = 1.038533918554399
$ws.Range("J22").Value = 1.039396214693293
$ws.Range("K22").Value = 1.041347821744816
$ws.Range("L22").Value = 1.044939336955539
$ws.Range("M22").Value = 1.053877538865288
$ws.Range("N22").Value = 1.040872276538569
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032885740035506
$ws.Range("D23").Value = 1.037845208175487
$ws.Range("E23").Value = 1.041559299459742
$ws.Range("F23").Value = 1.050570018374571
$ws.Range("I23").Value = 1.038606073237308
$ws.Range("J23").Value = 1.039609156313977
$ws.Range("K23").Value = 1.041479061510695
$ws.Range("L23").Value = 1.045179191202771
$ws.Range("M23").Value = 1.054156491815581
$ws.Range("N23").Value = 1.041085520560779
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034246095157549
$ws.Range("D24").Value = 1.038633917648343
$ws.Range("E24").Value = 1.042777559390406
$ws.Range("F24").Value = 1.051940398161287
$ws.Range("I24").Value = 1.038888047899519
$ws.Range("J24").Value = 1.040447465920503
$ws.Range("K24").Value = 1.041994813255075
$ws.Range("L24").Value = 1.046124138580186
$ws.Range("M24").Value = 1.05525575957299
$ws.Range("N24").Value = 1.041925020663102
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035827554373592
$ws.Range("D25").Value = 1.039550559618257
$ws.Range("E25").Value = 1.044195652042526
$ws.Range("F25").Value = 1.053535721608852
$ws.Range("I25").Value = 1.039210906618937
$ws.Range("J25").Value = 1.041420193662137
$ws.Range("K25").Value = 1.042591348422058
$ws.Range("L25").Value = 1.047222027848594
$ws.Range("M25").Value = 1.056533536237217
$ws.Range("N25").Value = 1.042899129789703
